# Amélioration du formulaire de contact et envoi du mail
#
# 1. Swap the tab order of "CU5 s'inscrire à la newsletter" and
#    "CU6 contacter l'association" (move the contact sheet one position to
#    the left, in front of the newsletter sheet), then rename both tabs so
#    the CU5/CU6 numbering again matches their (now swapped) position.
# 2. On the (now 5th, renamed) "contacter l'association" sheet, reword the
#    4th step of the nominal scenario.
# 3. Make the "contacter l'association" sheet the active sheet/selection.

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("CU6 contacter l'association")
$newsletter = $wb.Worksheets.Item("CU5 s'inscrire à la newsletter")

# Move "CU6 contacter l'association" so it sits right before the newsletter
# sheet - this swaps the two sheets' tab order (positions 5 and 6).
$contact.Move($newsletter)

# Re-acquire the sheets by their (now updated) tab position, since the
# worksheet handles above track position rather than the original sheet.
$pos5 = $wb.Worksheets.Item(5)
$pos6 = $wb.Worksheets.Item(6)

# Fix up the "CU5"/"CU6" numbering prefixes to match the new order.
$pos5.Name = "CU5 contacter l'association"
$pos6.Name = "CU6 s'inscrire à la newsletter"

# Reword step 4 of the nominal scenario on the contact sheet.
$pos5.Range("B18").Value = "4. Le site enregistre le message de l'internaute non membre ou du membre dans la base de données"

# Make the contact sheet active with the new selection.
$pos5.Activate()
$pos5.Range("E27").Select()
